$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# Header InlineShapes (BTec logo): rename from image1.jpg -> image2.jpg
$sec.Headers.Item(1).Range.InlineShapes.Item(1).Name = "image2.jpg"
$sec.Headers.Item(2).Range.InlineShapes.Item(1).Name = "image2.jpg"

# Footer InlineShapes (Pearson logo): rename from image2.png -> image1.png
$sec.Footers.Item(1).Range.InlineShapes.Item(1).Name = "image1.png"
$sec.Footers.Item(2).Range.InlineShapes.Item(1).Name = "image1.png"
